$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unmerge the ranges that are being resized (the bare single-cell
#     "merges" F1 / F2 / L2 already present in the source file cannot be
#     toggled back off through COM, so we only touch the real multi-cell
#     ones here) ---
$ws.Range("G1:L1").UnMerge()
$ws.Range("G2:K2").UnMerge()

# --- Row 1: drop the "2025" column header (G1) ---
$ws.Range("G1").Clear()

# --- Row 2: rename December -> Month 1, drop January / February cells ---
$ws.Range("F2").Value = "Month 1"
$ws.Range("G2").Clear()
$ws.Range("L2").Clear()

# --- Row 3: rename week labels, drop the extra week columns ---
$ws.Range("F3").Value = "Week 1"
$ws.Range("G3").Value = "Week 2"
$ws.Range("H3:L3").Clear()

# --- Row 4: Task 1 / ML1 - T1, drop start/end dates ---
$ws.Range("B4").Value = "Task 1"
$ws.Range("C4").Value = "ML1 - T1"
$ws.Range("D4").Clear()
$ws.Range("E4").Clear()

# --- Row 5: clear task number, drop start/end dates and highlight cell ---
$ws.Range("B5").ClearContents()
$ws.Range("D5").Clear()
$ws.Range("E5").Clear()
$ws.Range("G5").Clear()

# --- Row 6: Task 2 / ML2 - X1, move highlight from H6 to F6 ---
$ws.Range("B6").Value = "Task 2"
$ws.Range("C6").Value = "ML2 - X1"
$ws.Range("D6").Clear()
$ws.Range("E6").Clear()
$ws.Range("F4").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("H6").Clear()

# --- Row 7: Task 3 / ML2 - X2, move highlight from I7 to G7 ---
$ws.Range("B7").Value = "Task 3"
$ws.Range("C7").Value = "ML2 - X2"
$ws.Range("D7").Clear()
$ws.Range("E7").Clear()
$ws.Range("F4").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("I7").Clear()

# --- Row 8: Task 4 / ML2 - X3, move highlight from J8 to H8 ---
$ws.Range("B8").Value = "Task 4"
$ws.Range("C8").Value = "ML2 - X3"
$ws.Range("D8").Clear()
$ws.Range("E8").Clear()
$ws.Range("F4").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("J8").Clear()

$excel.CutCopyMode = $false

# --- Remove now-unused columns J:L ---
$ws.Range("J1:L8").EntireColumn.Delete()

# --- Re-merge the new ranges ---
$ws.Range("F1:G1").Merge()
$ws.Range("F2:I2").Merge()
